$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with freshly scraped values.
# D-column values are numeric-looking text (e.g. "1.001", "28.201.94") that must
# stay as plain text (matching the source data format), so we force text entry
# with a leading apostrophe and then reset the cell style so no extra
# number-format/quote-prefix styling is introduced.
$ws.Range("D2").Value = "'28.201.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("D3").Value = "'1.808.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.75%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'312.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.00%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("D7").Value = "'0.5143"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.19%  "
$ws.Range("D8").Value = "'0.3954"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.34%  "
$ws.Range("D9").Value = "'0.07812"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.34%  "
$ws.Range("E10").Value = "  -0.80%  "
$ws.Range("D11").Value = "'40.92"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.33%  "
$ws.Range("D12").Value = "'6.375"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("D13").Value = "'1.002"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("D14").Value = "'20.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.34%  "
$ws.Range("D15").Value = "'7.348"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.98%  "
$ws.Range("D16").Value = "'1.802.91"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.27%  "
$ws.Range("D17").Value = "'92.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.68%  "
$ws.Range("E18").Value = "  -2.22%  "
$ws.Range("D19").Value = "'0.06581"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.88%  "
$ws.Range("D20").Value = "'1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("E21").Value = "  -2.10%  "
$ws.Range("E22").Value = "  -0.38%  "
$ws.Range("D23").Value = "'28.251.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.68%  "
$ws.Range("E24").Value = "  -2.31%  "
$ws.Range("D25").Value = "'2.213"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.43%  "
$ws.Range("D26").Value = "'160.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.87%  "
$ws.Range("D27").Value = "'2.462"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.92%  "
$ws.Range("D28").Value = "'20.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.85%  "
$ws.Range("D29").Value = "'2.015.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.56%  "
$ws.Range("D30").Value = "'128.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.57%  "
$ws.Range("E31").Value = "  -0.75%  "
$ws.Range("D32").Value = "'1.061"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.05%  "
$ws.Range("D33").Value = "'3.659"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.48%  "
$ws.Range("D34").Value = "'5.571"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.95%  "
$ws.Range("D35").Value = "'0.07149"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.24%  "
$ws.Range("D36").Value = "'9.192"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.22%  "
$ws.Range("D37").Value = "'0.02350"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("D38").Value = "'0.2176"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.08%  "
$ws.Range("D39").Value = "'5.042"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.27%  "
$ws.Range("E40").Value = "  -7.00%  "
$ws.Range("D41").Value = "'0.6174"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.57%  "
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("D43").Value = "'1.151"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.79%  "
$ws.Range("D44").Value = "'13.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.55%  "
$ws.Range("D45").Value = "'0.5967"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.86%  "
$ws.Range("D46").Value = "'1.305"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.55%  "
$ws.Range("D47").Value = "'3.738"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.38%  "
$ws.Range("D48").Value = "'125.15"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.65%  "
$ws.Range("D49").Value = "'1.212"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.83%  "
$ws.Range("D50").Value = "'1.918"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.44%  "
$ws.Range("D51").Value = "'0.06797"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.39%  "
